# Fixed naive component forecaster bug - Presentation state 11.02.
# Shifts each data row's values one column to the right (B->C, C->D, ... J->K),
# dropping any value that falls past column K, and inserts a new "latest"
# value into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value to insert into column B for each data row (2-20).
$newB = @{
    2  = 0.3719860057927588
    3  = -2.702915518772638
    4  = -0.2307826431404359
    5  = -0.5654386276933741
    6  = -0.6603092772102132
    7  = -0.15162438770796
    8  = -0.2053460154962278
    9  = 0.6162032393936197
    10 = 1.652643173475852
    11 = 0.3110387314724781
    12 = 0.2388379152847414
    13 = 0.6508000635779043
    14 = 0.2387740594105157
    15 = 0.3465902496671606
    16 = 0.00230005330798793
    17 = -0.1902738424076751
    18 = -0.3325070745318338
    19 = 0.1656141382254278
    20 = -0.09587373626955231
}

$lastCol = 11  # column K

for ($r = 2; $r -le 20; $r++) {
    # Read existing values across the row (columns B..K), from left to right.
    $existing = @()
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $existing += $val
        }
    }

    # Shift values one column to the right, starting from the rightmost
    # occupied cell so we don't clobber values before they are read.
    $n = $existing.Count
    for ($i = $n - 1; $i -ge 0; $i--) {
        $destCol = 2 + $i + 1
        if ($destCol -le $lastCol) {
            $ws.Cells.Item($r, $destCol).Value = $existing[$i]
        }
    }

    # Insert the new value into column B.
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}
